# 20150603 - swap TIM1/TIM4 functions on AQ32 target
# (to better align with AQ32 pin usage.)
#
# This script updates the "Pins" worksheet of the AQ32/TauLabs hardware
# definition workbook:
#   - Renames the "AQ32 Use" / "TauLabs Use" column headers to
#     "AQ32 Definition" / "TauLabs Definition" (all three header blocks).
#   - Swaps the TIM1/TIM4 PWM-output and PPM/RangeFinder function names
#     between the two timers so the AQ32 target pin usage documentation
#     matches the firmware change.
#   - Updates the saved cell selection to K1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pins")

# ---------------------------------------------------------------------
# Header renames: "AQ32 Use" -> "AQ32 Definition"
#                 "TauLabs Use" -> "TauLabs Definition"
# (three repeated header blocks: B1/C1, F1/G1, J1/K1)
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "AQ32 Definition"
$ws.Range("C1").Value = "TauLabs Definition"
$ws.Range("F1").Value = "AQ32 Definition"
$ws.Range("G1").Value = "TauLabs Definition"
$ws.Range("J1").Value = "AQ32 Definition"
$ws.Range("K1").Value = "TauLabs Definition"

# ---------------------------------------------------------------------
# TIM4 column (F/G, rows 31-34): PE9-PE? area uses RX1..RX4 (TIM4 CHx)
# pin names in column F already; column G previously held
# "PWM Out 9-12 (TIM4 CHx)" -- now TIM4 is used for PPM/RangeFinder.
# ---------------------------------------------------------------------
$ws.Range("G31").Value = "PPM (TIM4 CH1)"
$ws.Range("G32").Value = "RangeFinder (TIM4 CH2)"
$ws.Range("G33").Clear()
$ws.Range("G34").Clear()

# ---------------------------------------------------------------------
# TIM1 column (J/K, rows 11-16): column J holds RX5-RX8 (TIM1 CHx) pin
# names already; column K previously only had entries for rows 15/16
# ("RangeFinder (TIM1 CH3)" / "Serial PPM (TIM1 CH4)"). TIM1 now
# provides PWM Out 9-12, so all four rows get "PWM Out N (TIM1 CHx)".
# ---------------------------------------------------------------------
$ws.Range("K11").Value = "PWM Out 9 (TIM1 CH1)"
$ws.Range("K13").Value = "PWM Out 10 (TIM1 CH2)"
$ws.Range("K15").Value = "PWM Out 11 (TIM1 CH3)"
$ws.Range("K16").Value = "PWM Out 12 (TIM1 CH4)"

# ---------------------------------------------------------------------
# Update saved selection to K1
# ---------------------------------------------------------------------
$ws.Range("K1").Select()
